$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-16 Tuesday" "2024-01-17 Wednesday"

Replace-Text "593×8=4744" "479×6=2874"
Replace-Text "844×4=3376" "279×3=837"
Replace-Text "257×8=2056" "551×8=4408"
Replace-Text "804×9=7236" "951×5=4755"
Replace-Text "490×4=1960" "963×6=5778"

Replace-Text "975×6=5850" "707×2=1414"
Replace-Text "185×2=370" "391×9=3519"
Replace-Text "568×3=1704" "363×9=3267"
Replace-Text "854×3=2562" "719×6=4314"
Replace-Text "972×8=7776" "352×7=2464"

Replace-Text "596×5=2980" "570×2=1140"
Replace-Text "121×9=1089" "681×8=5448"
Replace-Text "308×6=1848" "461×8=3688"
Replace-Text "843×3=2529" "320×7=2240"
Replace-Text "592×6=3552" "306×7=2142"

Replace-Text "141×7=987" "596×8=4768"
Replace-Text "928×7=6496" "293×6=1758"
Replace-Text "707×9=6363" "999×5=4995"
Replace-Text "225×9=2025" "448×4=1792"
Replace-Text "935×4=3740" "655×2=1310"

Replace-Text "943×3=2829" "909×8=7272"
Replace-Text "506×5=2530" "194×6=1164"
Replace-Text "263×2=526" "165×8=1320"
Replace-Text "635×2=1270" "469×3=1407"
Replace-Text "787×2=1574" "435×7=3045"
